$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.436801249702838
$ws.Range("C2").Value = 1.700317937929734
$ws.Range("D2").Value = 1.457347192194637
$ws.Range("E2").Value = 1.346982622402224
$ws.Range("B3").Value = 2.442838613609139
$ws.Range("C3").Value = 1.703193982803089
$ws.Range("D3").Value = 1.458595287136383
$ws.Range("E3").Value = 1.347828200706768
$ws.Range("B4").Value = 2.416199608207933
$ws.Range("C4").Value = 1.689607354824597
$ws.Range("D4").Value = 1.451431135937862
$ws.Range("E4").Value = 1.342718102337791
$ws.Range("B5").Value = 2.444040844258239
$ws.Range("C5").Value = 1.705018684914926
$ws.Range("D5").Value = 1.446166112413758
$ws.Range("E5").Value = 1.349992015156981
$ws.Range("B6").Value = 2.451383197835503
$ws.Range("C6").Value = 1.709059898859096
$ws.Range("D6").Value = 1.448580771627436
$ws.Range("E6").Value = 1.351843737668598
$ws.Range("B7").Value = 2.44341887967167
$ws.Range("C7").Value = 1.705004603960501
$ws.Range("D7").Value = 1.460838994595111
$ws.Range("E7").Value = 1.349955890124776
$ws.Range("B8").Value = 2.426137845639802
$ws.Range("C8").Value = 1.691113062267858
$ws.Range("D8").Value = 1.463289021859106
$ws.Range("E8").Value = 1.340676964773855
$ws.Range("B9").Value = 2.437293282939875
$ws.Range("C9").Value = 1.699213878649833
$ws.Range("D9").Value = 1.455487230471159
$ws.Range("E9").Value = 1.345344167322881
$ws.Range("B10").Value = 2.354271173204167
$ws.Range("C10").Value = 1.667191889118646
$ws.Range("D10").Value = 1.444906915254569
$ws.Range("E10").Value = 1.339021721330763
$ws.Range("B11").Value = 2.358096505327883
$ws.Range("C11").Value = 1.671843872822301
$ws.Range("D11").Value = 1.447518338027816
$ws.Range("E11").Value = 1.340293742975246
$ws.Range("B12").Value = 2.233541039399723
$ws.Range("C12").Value = 1.617368168486356
$ws.Range("D12").Value = 1.424266206036189
$ws.Range("E12").Value = 1.324371252076056
$ws.Range("B13").Value = 2.342009859240121
$ws.Range("C13").Value = 1.659963471116276
$ws.Range("D13").Value = 1.440926698548765
$ws.Range("E13").Value = 1.336020726629191
